$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected error copy-on-slice: fix Historic / Relative Change values in rows 45-48
$ws.Range("I45").Value = 17737.58
$ws.Range("J45").Value = 0.1299985680120963

$ws.Range("I46").Value = 75.01000000000001
$ws.Range("J46").Value = 0.2497000399946674

$ws.Range("I47").Value = 75.01000000000001
$ws.Range("J47").Value = 0.08732169044127448

$ws.Range("I48").Value = 75.01000000000001
$ws.Range("J48").Value = 0.7416344487401679

# Remove the now-obsolete UBA (2020) GreenSupreme meat-consumption row
$ws.Rows("49").Delete()
